<#
  Re-point the three "TestReport" data-bound content controls away from the
  custom XML part and onto the (empty) document Title core property - the
  same net effect Word produces when the bound custom XML data source is
  removed (e.g. via the Document Inspector "Remove All" on Custom XML
  Data): the controls fall back to a document-property binding that
  resolves to nothing, so they render their placeholder text.

  Also updates the literal paragraph text runs that sit in front of each
  content control ("TEST Name", "Details" -> "Client Name") and collapses
  the double space after "TEST Name".
#>

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>'
$pkgFooter = '</pkg:xmlData></pkg:part></pkg:package>'
$docNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Set-ParagraphXml($paraIndex, $innerParaXml) {
    $wrapped = $pkgHeader + '<w:document ' + $docNs + '><w:body>' + $innerParaXml + '</w:body></w:document>' + $pkgFooter
    $d.Paragraphs.Item($paraIndex).Range.InsertXML($wrapped)
}

# --- Paragraph 1: "TEST Name" content control -----------------------------
$p1 = '<w:p w14:paraId="39AE44E4" w14:textId="106BD432" w:rsidR="004A40D3" w:rsidRDefault="00117D17">' + `
    '<w:r><w:t xml:space="preserve">TEST Name </w:t></w:r>' + `
    '<w:sdt><w:sdtPr>' + `
        '<w:tag w:val="test_name"/>' + `
        '<w:id w:val="-520702687"/>' + `
        '<w:placeholder><w:docPart w:val="DefaultPlaceholder_-1854013440"/></w:placeholder>' + `
        '<w:showingPlcHdr/>' + `
        '<w:dataBinding w:prefixMappings="xmlns:ns0=''http://purl.org/dc/elements/1.1/'' xmlns:ns1=''http://schemas.openxmlformats.org/package/2006/metadata/core-properties'' " w:xpath="/ns1:coreProperties[1]/ns0:title[1]" w:storeItemID="{6C3C8BC8-F283-45AE-878A-BAB7291924A1}"/>' + `
        '<w:text/>' + `
    '</w:sdtPr><w:sdtContent>' + `
        '<w:r><w:rPr><w:rStyle w:val="PlaceholderText"/></w:rPr><w:t>Click or tap here to enter text.</w:t></w:r>' + `
    '</w:sdtContent></w:sdt>' + `
'</w:p>'

Set-ParagraphXml 1 $p1

# --- Paragraph 2: "TEST Description" content control ----------------------
$p2 = '<w:p w14:paraId="182AFC4C" w14:textId="19AE06B1" w:rsidR="00117D17" w:rsidRDefault="00117D17">' + `
    '<w:r><w:t xml:space="preserve">TEST Description </w:t></w:r>' + `
    '<w:sdt><w:sdtPr>' + `
        '<w:tag w:val="test_description"/>' + `
        '<w:id w:val="-888877327"/>' + `
        '<w:placeholder><w:docPart w:val="DefaultPlaceholder_-1854013440"/></w:placeholder>' + `
        '<w:showingPlcHdr/>' + `
        '<w:dataBinding w:prefixMappings="xmlns:ns0=''http://purl.org/dc/elements/1.1/'' xmlns:ns1=''http://schemas.openxmlformats.org/package/2006/metadata/core-properties'' " w:xpath="/ns1:coreProperties[1]/ns0:title[1]" w:storeItemID="{6C3C8BC8-F283-45AE-878A-BAB7291924A1}"/>' + `
        '<w:text/>' + `
    '</w:sdtPr><w:sdtContent>' + `
        '<w:r><w:rPr><w:rStyle w:val="PlaceholderText"/></w:rPr><w:t>Click or tap here to enter text.</w:t></w:r>' + `
    '</w:sdtContent></w:sdt>' + `
'</w:p>'

Set-ParagraphXml 2 $p2

# --- Paragraph 3: "Details" -> "Client Name" content control --------------
$p3 = '<w:p w14:paraId="2DD7DB27" w14:textId="7DA7B17C" w:rsidR="005820FD" w:rsidRPr="004A40D3" w:rsidRDefault="005820FD">' + `
    '<w:r><w:t xml:space="preserve">Client Name </w:t></w:r>' + `
    '<w:sdt><w:sdtPr>' + `
        '<w:tag w:val="settings_client_name"/>' + `
        '<w:id w:val="960917987"/>' + `
        '<w:placeholder><w:docPart w:val="DefaultPlaceholder_-1854013440"/></w:placeholder>' + `
        '<w:showingPlcHdr/>' + `
        '<w:dataBinding w:prefixMappings="xmlns:ns0=''http://purl.org/dc/elements/1.1/'' xmlns:ns1=''http://schemas.openxmlformats.org/package/2006/metadata/core-properties'' " w:xpath="/ns1:coreProperties[1]/ns0:title[1]" w:storeItemID="{6C3C8BC8-F283-45AE-878A-BAB7291924A1}"/>' + `
        '<w:text/>' + `
    '</w:sdtPr><w:sdtContent>' + `
        '<w:r><w:rPr><w:rStyle w:val="PlaceholderText"/></w:rPr><w:t>Click or tap here to enter text.</w:t></w:r>' + `
    '</w:sdtContent></w:sdt>' + `
'</w:p>'

Set-ParagraphXml 3 $p3

Write-Host "Done. Paragraphs:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host ("Para {0}: [{1}]" -f $i, $d.Paragraphs.Item($i).Range.Text)
}
